$d = $word.ActiveDocument
$RSQUOTE = [char]8217   # U+2019 RIGHT SINGLE QUOTATION MARK, used as the apostrophe in "X's"

# ---------------------------------------------------------------------------
# Edit 1: the stray "_GoBack" bookmark currently sits by itself in an
# otherwise empty paragraph (right before "The school record..."). It gets
# removed from there (it reappears elsewhere - see edit 3 below), leaving
# that paragraph as a plain empty paragraph.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# Edit 2: "Do you think Usain Bolt would have a similar shape graph to X's?"
# gets split into three separate runs:
#   "Do you think Usai" | "n Bolt" | " would have a similar shape graph to X's?"
# Dropping a zero-length bookmark at each split point and immediately
# deleting it again forces the run boundary to be created right there,
# without leaving behind any leftover run formatting or bookmark markup.
# ---------------------------------------------------------------------------
$fullText = $d.Content.Text
$needle = "Do you think Usain Bolt would have a similar shape graph to X" + $RSQUOTE + "s?"
$idx = $fullText.IndexOf($needle)
if ($idx -lt 0) {
    throw "Could not locate the 'Usain Bolt' sentence"
}

$prefix1 = "Do you think Usai"
$prefix2 = "Do you think Usain Bolt"
$split1 = $idx + $prefix1.Length
$split2 = $idx + $prefix2.Length

$d.Bookmarks.Add("__split1", $d.Range($split1, $split1))
$d.Bookmarks.Add("__split2", $d.Range($split2, $split2))
$d.Bookmarks.Item("__split1").Delete()
$d.Bookmarks.Item("__split2").Delete()

# ---------------------------------------------------------------------------
# Edit 3: "Do you think you would be able to use this data to help improve
# X's performance in the 100m sprint?" gets split right after "100m ", and
# the "_GoBack" bookmark removed in edit 1 is re-inserted exactly at that
# new split point (between the two resulting runs).
# ---------------------------------------------------------------------------
$fullText2 = $d.Content.Text
$needle2 = "Do you think you would be able to use this data to help improve X"
$idx2 = $fullText2.IndexOf($needle2)
if ($idx2 -lt 0) {
    throw "Could not locate the '100m sprint' sentence"
}

$prefix3 = "Do you think you would be able to use this data to help improve X" + $RSQUOTE + "s performance in the 100m "
$splitPos = $idx2 + $prefix3.Length

$d.Bookmarks.Add("_GoBack", $d.Range($splitPos, $splitPos))
